# "Einheiten und Combobox Längenbeschränkung"
#
# The underlying edit captured by the diff is a change to the input
# "Durchmesser" (diameter) value in cell D6 of Tabelle1, from 10 to 2.
# Every dependent formula in column H ("Vorschub pro Zahn in [mm]") is
# driven off of $D$6 (e.g. "=(0.025*$D$6)/2"), so updating D6 ripples
# through and recalculates all of those cached values automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D6").Value = 2
